# Generate Report for Handoff
# Update Priority ("low" -> "ht") and Latest Handoff Datetime for rows 4-7
# on both the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

for ($row = 4; $row -le 7; $row++) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-09-01 02:38:27"

    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-09-01 02:38:32"
}
